# Weekly update for "Hortaliza, Femacal de La Calera - Berenjena":
# a new observation is inserted as row 279 (pushing the existing rows
# 279-342 down to 280-343), growing the used range from A1:R342 to A1:R343.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 279; rows 279..342 shift down to 280..343.
$ws.Rows("279:279").Insert()

# Populate the newly inserted row with the new data point.
$ws.Cells.Item(279, 1).Value  = 3
$ws.Cells.Item(279, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(279, 3).Value  = "Coquimbo"
$ws.Cells.Item(279, 4).Value  = 44798
$ws.Cells.Item(279, 5).Value  = 5
$ws.Cells.Item(279, 6).Value  = 100112001
$ws.Cells.Item(279, 7).Value  = "Berenjena"
$ws.Cells.Item(279, 8).Value  = "Sin especificar"
$ws.Cells.Item(279, 9).Value  = "Primera"
$ws.Cells.Item(279, 10).Value = 105
$ws.Cells.Item(279, 11).Value = 8500
$ws.Cells.Item(279, 12).Value = 9000
$ws.Cells.Item(279, 13).Value = 8738
$ws.Cells.Item(279, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(279, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(279, 16).Value = 146
$ws.Cells.Item(279, 17).Value = 60
$ws.Cells.Item(279, 18).Value = "Hortaliza"
